$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.293"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05772"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.474"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.338"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8098"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8798"
$ws.Range("D9").Style = "Normal"
$ws.Range("B10").Value = "'One"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.01036"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'9OneONEBestin24h"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'WazirX"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.1377"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'10WazirXWRX"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.07283"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11MandalaExchangeTokenMDX"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.03080"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitrueCoin"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.03059"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'13BitrueCoinBTR"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'BitMartToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.09310"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'14BitMartTokenBMX"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'MCDex"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'3.854"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'15MCDexMCB"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'BitForexToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.001542"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'16BitForexTokenBF"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'CoinExToken"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'0.04700"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'17CoinExTokenCET"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.006058"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001302"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004603"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00008810"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'21NitroExNTX"
$ws.Range("E22").Style = "Normal"
$ws.Range("D26").Value = "'0.1319"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = "'0.0002352"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.03764"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006327"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = "'CEJI"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.004005"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'41CEJICEJI"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'BKEXToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.1051"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'42BKEXTokenBKK"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007120"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005474"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.6007"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.001864"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'47BOLOBOLOWorstin24h"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
